$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Label" header in H1, matching the style of the other headers (B1:G1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the Label column (H2:H15) with the values from the diff
$labels = @(0, 0, 1, 1, 1, 1, 1, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
